$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196; this pushes existing rows 196:251 down to 197:252
$ws.Rows.Item(196).Insert()

# Populate the new row 196 with its data (same static columns as its neighbours,
# new date and price figures per the commit)
$ws.Range("A196").Value = 5
$ws.Range("B196").Value = "Macroferia Regional de Talca"
$ws.Range("C196").Value = "Maule"
$ws.Range("D196").Value = 44722
$ws.Range("E196").Value = 7
$ws.Range("F196").Value = 100112008
$ws.Range("G196").Value = "Coliflor"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 2000
$ws.Range("K196").Value = 1200
$ws.Range("L196").Value = 1200
$ws.Range("M196").Value = 1200
$ws.Range("N196").Value = "$/unidad"
$ws.Range("O196").Value = "Región del Maule"
$ws.Range("P196").Value = 1200
$ws.Range("Q196").Value = 1
$ws.Range("R196").Value = "Hortaliza"
